# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values for column G (rows 2-26), replacing the old Strike# counts.
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 4
    8  = 2
    9  = 0
    10 = 3
    11 = 1
    12 = 5
    13 = 7
    14 = 1
    15 = 0
    16 = 1
    17 = 3
    18 = 2
    19 = 5
    20 = 2
    21 = 2
    22 = 1
    23 = 6
    24 = 6
    25 = 3
    26 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
